# Auto-generated edit script: apply numeric corrections to Leve profit sheets
# per scheduled-runner recompute (commit: chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1049.5
$ws.Range("I12").Value = 399.33334
$ws.Range("K12").Value = 399.33334
$ws.Range("M12").Value = -229.33334
$ws.Range("H15").Value = 1862.3846
$ws.Range("I15").Value = 1862.3846
$ws.Range("K15").Value = 5587.1538
$ws.Range("M15").Value = -5418.1538
$ws.Range("H17").Value = 1054894.2
$ws.Range("I17").Value = 884.0
$ws.Range("J17").Value = 1713650.6
$ws.Range("K17").Value = 2652.0
$ws.Range("L17").Value = 5140951.800000001
$ws.Range("M17").Value = -2484.0
$ws.Range("N17").Value = -5141287.800000001
$ws.Range("H40").Value = 5559260.5
$ws.Range("I40").Value = 2776.6667
$ws.Range("K40").Value = 2776.6667
$ws.Range("M40").Value = -2601.6667
$ws.Range("H100").Value = 1200.3636
$ws.Range("I100").Value = 1098.5
$ws.Range("J100").Value = 1258.5714
$ws.Range("K100").Value = 1098.5
$ws.Range("L100").Value = 1258.5714
$ws.Range("M100").Value = -557.5
$ws.Range("N100").Value = -2340.5714
$ws.Range("H116").Value = 27787862.0
$ws.Range("J116").Value = 11345.625
$ws.Range("L116").Value = 11345.625
$ws.Range("N116").Value = -18229.625
$ws.Range("H132").Value = 1823.0217
$ws.Range("I132").Value = 1727.027
$ws.Range("J132").Value = 2217.6667
$ws.Range("K132").Value = 5181.081
$ws.Range("L132").Value = 6653.000100000001
$ws.Range("M132").Value = -2651.081
$ws.Range("N132").Value = -11713.0001
$ws.Range("H137").Value = 9874.967
$ws.Range("I137").Value = 8105.294
$ws.Range("K137").Value = 24315.882
$ws.Range("M137").Value = -21765.882
$ws.Range("H138").Value = 4981.8335
$ws.Range("J138").Value = 5269.52
$ws.Range("L138").Value = 15808.56
$ws.Range("N138").Value = -26088.56

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3654.41
$ws.Range("I32").Value = 3458.1545
$ws.Range("K32").Value = 3458.1545
$ws.Range("M32").Value = -3171.1545
$ws.Range("H45").Value = 3881.182
$ws.Range("I45").Value = 0.0
$ws.Range("J45").Value = 3881.182
$ws.Range("K45").Value = 0.0
$ws.Range("L45").Value = 3881.182
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -4635.182
$ws.Range("H74").Value = 25346.092
$ws.Range("I74").Value = 32777.625
$ws.Range("K74").Value = 32777.625
$ws.Range("M74").Value = -31903.625
$ws.Range("H77").Value = 25346.092
$ws.Range("I77").Value = 32777.625
$ws.Range("K77").Value = 163888.125
$ws.Range("M77").Value = -159520.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4018.6316
$ws.Range("I94").Value = 1941.75
$ws.Range("K94").Value = 1941.75
$ws.Range("M94").Value = -1490.75
$ws.Range("H141").Value = 59746.0
$ws.Range("J141").Value = 59709.855
$ws.Range("L141").Value = 59709.855
$ws.Range("N141").Value = -70069.85500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2694.5715
$ws.Range("I16").Value = 1474.8334
$ws.Range("K16").Value = 1474.8334
$ws.Range("M16").Value = -1187.8334
$ws.Range("H22").Value = 1388.875
$ws.Range("I22").Value = 1450.1428
$ws.Range("J22").Value = 960.0
$ws.Range("K22").Value = 1450.1428
$ws.Range("L22").Value = 960.0
$ws.Range("M22").Value = -1100.1428
$ws.Range("N22").Value = -1660.0
$ws.Range("H31").Value = 4951.75
$ws.Range("I31").Value = 1990.7949
$ws.Range("K31").Value = 1990.7949
$ws.Range("M31").Value = -1695.7949
$ws.Range("H34").Value = 4951.75
$ws.Range("I34").Value = 1990.7949
$ws.Range("K34").Value = 1990.7949
$ws.Range("M34").Value = -1788.7949
$ws.Range("H50").Value = 79067.2
$ws.Range("J50").Value = 79067.2
$ws.Range("L50").Value = 79067.2
$ws.Range("N50").Value = -80317.2
$ws.Range("H51").Value = 44300.0
$ws.Range("J51").Value = 44380.0
$ws.Range("L51").Value = 44380.0
$ws.Range("N51").Value = -45852.0
$ws.Range("H61").Value = 44300.0
$ws.Range("J61").Value = 44380.0
$ws.Range("L61").Value = 44380.0
$ws.Range("N61").Value = -45076.0
$ws.Range("H99").Value = 5904.0
$ws.Range("I99").Value = 4679.8
$ws.Range("K99").Value = 4679.8
$ws.Range("M99").Value = -3181.8
$ws.Range("H106").Value = 50591.0
$ws.Range("J106").Value = 50591.0
$ws.Range("L106").Value = 50591.0
$ws.Range("N106").Value = -53115.0
$ws.Range("H107").Value = 2445.4688
$ws.Range("I107").Value = 1891.4736
$ws.Range("K107").Value = 1891.4736
$ws.Range("M107").Value = 28.52639999999997
$ws.Range("H113").Value = 2694.5715
$ws.Range("I113").Value = 1474.8334
$ws.Range("K113").Value = 1474.8334
$ws.Range("M113").Value = 695.1666
$ws.Range("H122").Value = 4543.5293
$ws.Range("I122").Value = 3253.5
$ws.Range("K122").Value = 9760.5
$ws.Range("M122").Value = -7310.5
$ws.Range("H126").Value = 5904.0
$ws.Range("I126").Value = 4679.8
$ws.Range("K126").Value = 14039.4
$ws.Range("M126").Value = -11569.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 9749.444
$ws.Range("J39").Value = 14199.0
$ws.Range("L39").Value = 42597.0
$ws.Range("N39").Value = -43185.0
$ws.Range("H50").Value = 25642336.0
$ws.Range("I50").Value = 66667650.0
$ws.Range("J50").Value = 1514.375
$ws.Range("K50").Value = 200002950.0
$ws.Range("L50").Value = 4543.125
$ws.Range("M50").Value = -200002469.0
$ws.Range("N50").Value = -5505.125
$ws.Range("H53").Value = 25642336.0
$ws.Range("I53").Value = 66667650.0
$ws.Range("J53").Value = 1514.375
$ws.Range("K53").Value = 200002950.0
$ws.Range("L53").Value = 4543.125
$ws.Range("M53").Value = -200002469.0
$ws.Range("N53").Value = -5505.125
$ws.Range("H131").Value = 30696.0
$ws.Range("I131").Value = 2485.6667
$ws.Range("J131").Value = 33340.72
$ws.Range("K131").Value = 7457.000100000001
$ws.Range("L131").Value = 100022.16
$ws.Range("M131").Value = -2417.000100000001
$ws.Range("N131").Value = -110102.16
$ws.Range("H132").Value = 7950.378
$ws.Range("I132").Value = 5423.778
$ws.Range("J132").Value = 9634.777
$ws.Range("K132").Value = 48814.002
$ws.Range("L132").Value = 86712.993
$ws.Range("M132").Value = -46284.002
$ws.Range("N132").Value = -91772.993
$ws.Range("H137").Value = 253636.12
$ws.Range("I137").Value = 171342.83
$ws.Range("J137").Value = 500516.0
$ws.Range("K137").Value = 514028.49
$ws.Range("L137").Value = 1501548.0
$ws.Range("M137").Value = -508928.49
$ws.Range("N137").Value = -1511748.0
$ws.Range("H141").Value = 4604.278
$ws.Range("I141").Value = 4604.278
$ws.Range("K141").Value = 13812.834
$ws.Range("M141").Value = -8632.834

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 6764.75
$ws.Range("I55").Value = 4686.3335
$ws.Range("K55").Value = 4686.3335
$ws.Range("M55").Value = -4359.3335
$ws.Range("H80").Value = 3747.5
$ws.Range("I80").Value = 3500.0
$ws.Range("K80").Value = 3500.0
$ws.Range("M80").Value = -2502.0
$ws.Range("H83").Value = 3747.5
$ws.Range("I83").Value = 3500.0
$ws.Range("K83").Value = 17500.0
$ws.Range("M83").Value = -12508.0
$ws.Range("H122").Value = 6586701.0
$ws.Range("I122").Value = 6586701.0
$ws.Range("J122").Value = 0.0
$ws.Range("K122").Value = 19760103.0
$ws.Range("L122").Value = 0.0
$ws.Range("M122").Value = -19757653.0
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 30000.0
$ws.Range("J123").Value = 30000.0
$ws.Range("L123").Value = 30000.0
$ws.Range("N123").Value = -34900.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1553.1177
$ws.Range("I46").Value = 1201.1666
$ws.Range("K46").Value = 1201.1666
$ws.Range("M46").Value = -1013.1666
$ws.Range("H68").Value = 6124.55
$ws.Range("I68").Value = 6124.25
$ws.Range("J68").Value = 6124.625
$ws.Range("K68").Value = 6124.25
$ws.Range("L68").Value = 6124.625
$ws.Range("M68").Value = -5375.25
$ws.Range("N68").Value = -7622.625
$ws.Range("H71").Value = 6124.55
$ws.Range("I71").Value = 6124.25
$ws.Range("J71").Value = 6124.625
$ws.Range("K71").Value = 30621.25
$ws.Range("L71").Value = 30623.125
$ws.Range("M71").Value = -26877.25
$ws.Range("N71").Value = -38111.125
$ws.Range("H93").Value = 3841.8572
$ws.Range("I93").Value = 4779.6
$ws.Range("J93").Value = 1497.5
$ws.Range("K93").Value = 4779.6
$ws.Range("L93").Value = 1497.5
$ws.Range("M93").Value = -3531.6
$ws.Range("N93").Value = -3993.5
$ws.Range("H136").Value = 10884.207
$ws.Range("I136").Value = 3897.3076
$ws.Range("K136").Value = 11691.9228
$ws.Range("M136").Value = -9141.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1364775.8
$ws.Range("J81").Value = 4340.3335
$ws.Range("L81").Value = 8680.667
$ws.Range("N81").Value = -10802.667
$ws.Range("H84").Value = 1364775.8
$ws.Range("J84").Value = 4340.3335
$ws.Range("L84").Value = 43403.335
$ws.Range("N84").Value = -54011.335
$ws.Range("H104").Value = 21315.625
$ws.Range("J104").Value = 21315.625
$ws.Range("L104").Value = 21315.625
$ws.Range("N104").Value = -28303.625
